$wb = $excel.ActiveWorkbook

# --- Update Hoja1!A1 conversion text with new rates ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 3.36 = 12969.8 pesos`n✅ 12969.8 pesos = 3.34 = 956.02 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

$ws1.Range("A1").Value = $newText

# --- Update tasas sheet N10/O10/N12/O12 values ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 298
$ws2.Range("O10").Value = 3865
$ws2.Range("N12").Value = 3880
$ws2.Range("O12").Value = 286
